$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class 0.0)
$ws.Range("B2").Value = 0.9504814305364512
$ws.Range("C2").Value = 0.8680904522613065
$ws.Range("D2").Value = 0.9074195666447801
$ws.Range("E2").Value = 796

# Row 3 (class 1.0)
$ws.Range("B3").Value = 0.6125461254612546
$ws.Range("C3").Value = 0.8217821782178217
$ws.Range("D3").Value = 0.7019027484143764

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.8587174348697395
$ws.Range("C4").Value = 0.8587174348697395
$ws.Range("D4").Value = 0.8587174348697395
$ws.Range("E4").Value = 0.8587174348697395

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.7815137779988529
$ws.Range("C5").Value = 0.8449363152395641
$ws.Range("D5").Value = 0.8046611575295782
$ws.Range("E5").Value = 998

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.8820816994490868
$ws.Range("C6").Value = 0.8587174348697395
$ws.Range("D6").Value = 0.8658219741773036
$ws.Range("E6").Value = 998
